$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 221, shifting existing rows 221:309 down to 222:310
$ws.Rows(221).Insert()

# Populate the new row 221 with the new weekly data point
$ws.Range("A221").Value = 6
$ws.Range("B221").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C221").Value = "Metropolitana"
$ws.Range("D221").Value = 45009
$ws.Range("E221").Value = 13
$ws.Range("F221").Value = 100112029
$ws.Range("G221").Value = "Orégano"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 48
$ws.Range("K221").Value = 17000
$ws.Range("L221").Value = 18000
$ws.Range("M221").Value = 17458
$ws.Range("N221").Value = '$/docena de atados'
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 5819
$ws.Range("Q221").Value = 3
$ws.Range("R221").Value = "Hortaliza"
